$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 10444
$ws.Range("F6").Value = 618
$ws.Range("F7").Value = 87
$ws.Range("F8").Value = 1739
$ws.Range("F11").Value = 252
$ws.Range("F12").Value = 553
$ws.Range("F13").Value = 1182
$ws.Range("F14").Value = 151
$ws.Range("F16").Value = 1032
$ws.Range("F18").Value = 124
$ws.Range("F19").Value = 434
$ws.Range("F20").Value = 434
$ws.Range("F21").Value = 23
$ws.Range("F22").Value = 368
$ws.Range("F23").Value = 67
$ws.Range("F24").Value = 1109
$ws.Range("F25").Value = 1126
$ws.Range("F26").Value = 1236
$ws.Range("F27").Value = 222
$ws.Range("F28").Value = 1427
$ws.Range("F29").Value = 727
$ws.Range("F31").Value = 34
$ws.Range("F32").Value = 100
$ws.Range("F33").Value = 709
$ws.Range("F37").Value = 802
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 820
$ws.Range("F40").Value = 535
$ws.Range("F41").Value = 1280
$ws.Range("F42").Value = 837
$ws.Range("F43").Value = 762
$ws.Range("F44").Value = 1395
$ws.Range("F46").Value = 731

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 121
$ws.Range("F19").Value = 1167
$ws.Range("F21").Value = 2299
$ws.Range("F22").Value = 1124
$ws.Range("F25").Value = 103
$ws.Range("F30").Value = 387
$ws.Range("F31").Value = 12
$ws.Range("F33").Value = 230
$ws.Range("F39").Value = 37
$ws.Range("F43").Value = 16
$ws.Range("F47").Value = 44

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 217
$ws.Range("F6").Value = 2587
$ws.Range("F7").Value = 4272
$ws.Range("F10").Value = 422
$ws.Range("F11").Value = 394
$ws.Range("F12").Value = 281
$ws.Range("F13").Value = 207
$ws.Range("F14").Value = 91

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 10444
$ws.Range("F5").Value = 217
$ws.Range("F6").Value = 4272
$ws.Range("F7").Value = 618
$ws.Range("F8").Value = 394
$ws.Range("F9").Value = 1739
$ws.Range("F10").Value = 252
$ws.Range("F13").Value = 207
$ws.Range("F14").Value = 91
$ws.Range("F15").Value = 121
$ws.Range("F16").Value = 1032
$ws.Range("F18").Value = 124
$ws.Range("F19").Value = 434
$ws.Range("F20").Value = 434
$ws.Range("F21").Value = 368
$ws.Range("F22").Value = 2299
$ws.Range("F23").Value = 2299
$ws.Range("F24").Value = 1124
$ws.Range("F25").Value = 1109
$ws.Range("F26").Value = 1126
$ws.Range("F27").Value = 1236
$ws.Range("F28").Value = 103
$ws.Range("F30").Value = 1427
$ws.Range("F31").Value = 727
$ws.Range("F32").Value = 387
$ws.Range("F33").Value = 709
$ws.Range("F34").Value = 12
$ws.Range("F37").Value = 802
$ws.Range("F38").Value = 230
$ws.Range("F39").Value = 820
$ws.Range("F40").Value = 535
$ws.Range("F41").Value = 837
$ws.Range("F43").Value = 762
$ws.Range("F44").Value = 1395
$ws.Range("F48").Value = 731
